$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.814.46"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.35%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.601.17"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.54%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "523.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.97%  "

$ws.Range("E6").Value = "  -1.33%  "

$ws.Range("E8").Value = "  +0.87%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.68"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.19%  "

$ws.Range("E10").Value = "  +1.38%  "

$ws.Range("E11").Value = "  -0.53%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.055.69"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.54%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "60.824.02"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.41%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.58"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.76%  "

$ws.Range("E16").Value = "  +0.57%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.602.57"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.35%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.76"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.52%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "355.01"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.27%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.55"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.31%  "

$ws.Range("E21").Value = "  +0.89%  "

$ws.Range("E22").Value = "  +0.25%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "61.11"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.14%  "

$ws.Range("E24").Value = "  +1.02%  "

$ws.Range("E25").Value = "  -0.51%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.716.89"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.72%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.996"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.32%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0843"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.81%  "

$ws.Range("E29").Value = "  -0.84%  "

$ws.Range("E30").Value = "  +0.02%  "

$ws.Range("E31").Value = "  +8.85%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.37"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.20%  "

$ws.Range("E33").Value = "  +2.04%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "150.23"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.72%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.17"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.88%  "

$ws.Range("E36").Value = "  +7.68%  "

$ws.Range("E37").Value = "  +0.08%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.903"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.30%  "

$ws.Range("E39").Value = "  +0.56%  "

$ws.Range("E40").Value = "  +0.50%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "297.80"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.18%  "

$ws.Range("E42").Value = "  +1.30%  "

$ws.Range("E43").Value = "  -0.04%  "

$ws.Range("E44").Value = "  -0.54%  "

$ws.Range("E45").Value = "  +0.12%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.57"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.22%  "

$ws.Range("E47").Value = "  +0.91%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0237"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.22%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.33"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.40%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.23"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.27%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.968.90"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.60%  "

